$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: B31 was an empty inline-string cell in the source feed; the bot run
# that produced this edit left it blank this time, so clear it out entirely.
$ws.Range("B31").ClearContents()

# Row 32: new article scraped by the newsbot
$ws.Range("A32").Value = "04/01/2026 23:04:18"
$ws.Range("B32").Value = "04/01 23:00"
$ws.Range("C32").Value = "Folha de S.Paulo - Mercado - Principal"
$ws.Range("D32").Value = "Instituições financeiras vão reforçar apoio a BC após TCU anunciar inspeção no caso Master"
$ws.Range("E32").Value = "https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/mercado/2026/01/instituicoes-financeiras-vao-reforcar-apoio-a-bc-apos-tcu-anunciar-inspecao-no-caso-master.shtml"
$ws.Range("F32").Value = 2
$ws.Range("G32").Value = "tcu"
$ws.Range("H32").Value = "ira brasileira decidiu reforçar seu apoio ao Banco Central após o presidente do TCU (Tribunal de Contas da Conta), Vital do Rêgo, determinar a inspeção in loco da "

# Row 33: new article scraped by the newsbot
$ws.Range("A33").Value = "04/01/2026 23:04:19"
$ws.Range("B33").Value = "04/01 23:00"
$ws.Range("C33").Value = "Folha de S.Paulo - Mercado - Principal"
$ws.Range("D33").Value = "Cotado para suceder Haddad, Dario Durigan é conhecido como 'CEO' do Ministério da Fazenda"
$ws.Range("E33").Value = "https://redir.folha.com.br/redir/online/mercado/rss091/*https://www1.folha.uol.com.br/mercado/2026/01/cotado-para-suceder-haddad-dario-durigan-e-conhecido-como-ceo-do-ministerio-da-fazenda.shtml"
$ws.Range("F33").Value = 2
$ws.Range("G33").Value = "haddad"
$ws.Range("H33").Value = "ogado Dario Durigan, 41, &lt;a href=&quot;https://www1.folha.uol.com.br/mercado/2023/06/<b>haddad</b>-troca-perfil-politico-por-gestor-com-novo-numero-2-da-fazenda.shtml&quot;&gt;assumiu a "
